$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value2 = 6362011
$ws.Range("F3").Value2 = "AE Altos"
$ws.Range("G3").Value2 = "CSA"
$ws.Range("H3").Value2 = 1
$ws.Range("I3").Value2 = 1
$ws.Range("J3").Value2 = "D"
$ws.Range("K3").Value2 = 2.8
$ws.Range("L3").Value2 = 2.9
$ws.Range("M3").Value2 = 2.4
$ws.Range("N3").Value2 = 2.6
$ws.Range("O3").Value2 = 2.875
$ws.Range("P3").Value2 = 2.625
$ws.Range("Q3").Value2 = 0
$ws.Range("R3").Value2 = 1.875
$ws.Range("S3").Value2 = 1.925
$ws.Range("T3").Value2 = 2
$ws.Range("U3").Value2 = 1.925
$ws.Range("V3").Value2 = 1.875
$ws.Range("W3").Value2 = -1
$ws.Range("X3").Value2 = 1.875
$ws.Range("Y3").Value2 = -1
$ws.Range("Z3").Value2 = 0
$ws.Range("AA3").Value2 = -0
$ws.Range("AB3").Value2 = 0
$ws.Range("AC3").Value2 = -0

# Row 4
$ws.Range("B4").Value2 = 6361796
$ws.Range("F4").Value2 = "Manaus"
$ws.Range("G4").Value2 = "Volta Redonda"
$ws.Range("H4").Value2 = 1
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = "H"
$ws.Range("K4").Value2 = 2.375
$ws.Range("L4").Value2 = 3
$ws.Range("M4").Value2 = 2.8
$ws.Range("N4").Value2 = 2.55
$ws.Range("O4").Value2 = 2.8
$ws.Range("P4").Value2 = 2.8
$ws.Range("Q4").Value2 = 0
$ws.Range("R4").Value2 = 1.775
$ws.Range("S4").Value2 = 2.025
$ws.Range("T4").Value2 = 2
$ws.Range("U4").Value2 = 1.875
$ws.Range("V4").Value2 = 1.925
$ws.Range("W4").Value2 = 1.55
$ws.Range("X4").Value2 = -1
$ws.Range("Y4").Value2 = -1
$ws.Range("Z4").Value2 = 0.7749999999999999
$ws.Range("AA4").Value2 = -1
$ws.Range("AB4").Value2 = -1
$ws.Range("AC4").Value2 = 0.925

# Row 7
$ws.Range("B7").Value2 = 6361795
$ws.Range("F7").Value2 = "Botafogo PB"
$ws.Range("G7").Value2 = "Floresta EC"
$ws.Range("H7").Value2 = 1
$ws.Range("I7").Value2 = 1
$ws.Range("J7").Value2 = "D"
$ws.Range("K7").Value2 = 1.533
$ws.Range("L7").Value2 = 3.5
$ws.Range("M7").Value2 = 5.5
$ws.Range("N7").Value2 = 1.5
$ws.Range("O7").Value2 = 3.6
$ws.Range("P7").Value2 = 5.5
$ws.Range("Q7").Value2 = -0.75
$ws.Range("R7").Value2 = 1.7
$ws.Range("S7").Value2 = 2.1
$ws.Range("T7").Value2 = 2
$ws.Range("U7").Value2 = 1.8
$ws.Range("V7").Value2 = 2
$ws.Range("W7").Value2 = -1
$ws.Range("X7").Value2 = 2.6
$ws.Range("Y7").Value2 = -1
$ws.Range("Z7").Value2 = -1
$ws.Range("AA7").Value2 = 1.1
$ws.Range("AB7").Value2 = 0
$ws.Range("AC7").Value2 = -0

# Row 8
$ws.Range("B8").Value2 = 6361626
$ws.Range("F8").Value2 = "Paysandu"
$ws.Range("G8").Value2 = "Sao Jose PA"
$ws.Range("H8").Value2 = 1
$ws.Range("I8").Value2 = 1
$ws.Range("J8").Value2 = "D"
$ws.Range("K8").Value2 = 1.85
$ws.Range("L8").Value2 = 3
$ws.Range("M8").Value2 = 4.2
$ws.Range("N8").Value2 = 1.8
$ws.Range("O8").Value2 = 3
$ws.Range("P8").Value2 = 4.5
$ws.Range("Q8").Value2 = -0.5
$ws.Range("R8").Value2 = 1.85
$ws.Range("S8").Value2 = 1.95
$ws.Range("T8").Value2 = 2
$ws.Range("U8").Value2 = 1.9
$ws.Range("V8").Value2 = 1.9
$ws.Range("W8").Value2 = -1
$ws.Range("X8").Value2 = 2
$ws.Range("Y8").Value2 = -1
$ws.Range("Z8").Value2 = -1
$ws.Range("AA8").Value2 = 0.95
$ws.Range("AB8").Value2 = 0
$ws.Range("AC8").Value2 = -0

# Row 24
$ws.Range("B24").Value2 = 6758548
$ws.Range("F24").Value2 = "Operario PR"
$ws.Range("G24").Value2 = "CSA"
$ws.Range("H24").Value2 = 1
$ws.Range("I24").Value2 = 0
$ws.Range("J24").Value2 = "H"
$ws.Range("K24").Value2 = 2
$ws.Range("L24").Value2 = 2.9
$ws.Range("M24").Value2 = 3.75
$ws.Range("N24").Value2 = 1.909
$ws.Range("O24").Value2 = 3
$ws.Range("P24").Value2 = 4.2
$ws.Range("Q24").Value2 = -0.5
$ws.Range("R24").Value2 = 1.95
$ws.Range("S24").Value2 = 1.85
$ws.Range("T24").Value2 = 2
$ws.Range("U24").Value2 = 1.85
$ws.Range("V24").Value2 = 1.95
$ws.Range("W24").Value2 = 0.909
$ws.Range("X24").Value2 = -1
$ws.Range("Y24").Value2 = -1
$ws.Range("Z24").Value2 = 0.95
$ws.Range("AA24").Value2 = -1
$ws.Range("AB24").Value2 = -1
$ws.Range("AC24").Value2 = 0.95

# Row 25
$ws.Range("B25").Value2 = 6362384
$ws.Range("F25").Value2 = "Aparecidense"
$ws.Range("G25").Value2 = "Clube Do Remo"
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 2
$ws.Range("J25").Value2 = "A"
$ws.Range("K25").Value2 = 2.4
$ws.Range("L25").Value2 = 2.8
$ws.Range("M25").Value2 = 3
$ws.Range("N25").Value2 = 2.4
$ws.Range("O25").Value2 = 2.75
$ws.Range("P25").Value2 = 3.1
$ws.Range("Q25").Value2 = -0.25
$ws.Range("R25").Value2 = 2.025
$ws.Range("S25").Value2 = 1.775
$ws.Range("T25").Value2 = 2
$ws.Range("U25").Value2 = 1.85
$ws.Range("V25").Value2 = 1.95
$ws.Range("W25").Value2 = -1
$ws.Range("X25").Value2 = -1
$ws.Range("Y25").Value2 = 2.1
$ws.Range("Z25").Value2 = -1
$ws.Range("AA25").Value2 = 0.7749999999999999
$ws.Range("AB25").Value2 = 0
$ws.Range("AC25").Value2 = -0

# Row 65
$ws.Range("B65").Value2 = 6361946
$ws.Range("F65").Value2 = "Pouso Alegre"
$ws.Range("G65").Value2 = "Floresta EC"
$ws.Range("H65").Value2 = 0
$ws.Range("I65").Value2 = 1
$ws.Range("J65").Value2 = "A"
$ws.Range("K65").Value2 = 2.2
$ws.Range("L65").Value2 = 2.875
$ws.Range("M65").Value2 = 3.2
$ws.Range("N65").Value2 = 1.85
$ws.Range("O65").Value2 = 3.1
$ws.Range("P65").Value2 = 3.8
$ws.Range("Q65").Value2 = -0.5
$ws.Range("R65").Value2 = 1.9
$ws.Range("S65").Value2 = 1.9
$ws.Range("T65").Value2 = 2
$ws.Range("U65").Value2 = 1.95
$ws.Range("V65").Value2 = 1.85
$ws.Range("W65").Value2 = -1
$ws.Range("X65").Value2 = -1
$ws.Range("Y65").Value2 = 2.8
$ws.Range("Z65").Value2 = -1
$ws.Range("AA65").Value2 = 0.8999999999999999
$ws.Range("AB65").Value2 = -1
$ws.Range("AC65").Value2 = 0.8500000000000001

# Row 66
$ws.Range("B66").Value2 = 6361488
$ws.Range("F66").Value2 = "America RN"
$ws.Range("G66").Value2 = "Manaus"
$ws.Range("H66").Value2 = 1
$ws.Range("I66").Value2 = 2
$ws.Range("J66").Value2 = "A"
$ws.Range("K66").Value2 = 1.533
$ws.Range("L66").Value2 = 3.4
$ws.Range("M66").Value2 = 6
$ws.Range("N66").Value2 = 1.5
$ws.Range("O66").Value2 = 3.5
$ws.Range("P66").Value2 = 6
$ws.Range("Q66").Value2 = -1
$ws.Range("R66").Value2 = 1.975
$ws.Range("S66").Value2 = 1.825
$ws.Range("T66").Value2 = 2.25
$ws.Range("U66").Value2 = 2.025
$ws.Range("V66").Value2 = 1.775
$ws.Range("W66").Value2 = -1
$ws.Range("X66").Value2 = -1
$ws.Range("Y66").Value2 = 5
$ws.Range("Z66").Value2 = -1
$ws.Range("AA66").Value2 = 0.825
$ws.Range("AB66").Value2 = 1.025
$ws.Range("AC66").Value2 = -1

# Row 67
$ws.Range("B67").Value2 = 6362457
$ws.Range("F67").Value2 = "Figueirense"
$ws.Range("G67").Value2 = "CSA"
$ws.Range("H67").Value2 = 1
$ws.Range("I67").Value2 = 1
$ws.Range("J67").Value2 = "D"
$ws.Range("K67").Value2 = 2.25
$ws.Range("L67").Value2 = 2.875
$ws.Range("M67").Value2 = 3.1
$ws.Range("N67").Value2 = 2.15
$ws.Range("O67").Value2 = 2.9
$ws.Range("P67").Value2 = 3.25
$ws.Range("Q67").Value2 = -0.25
$ws.Range("R67").Value2 = 1.9
$ws.Range("S67").Value2 = 1.9
$ws.Range("T67").Value2 = 1.75
$ws.Range("U67").Value2 = 1.75
$ws.Range("V67").Value2 = 2.05
$ws.Range("W67").Value2 = -1
$ws.Range("X67").Value2 = 1.9
$ws.Range("Y67").Value2 = -1
$ws.Range("Z67").Value2 = -0.5
$ws.Range("AA67").Value2 = 0.45
$ws.Range("AB67").Value2 = 0.375
$ws.Range("AC67").Value2 = -0.5

# Row 68
$ws.Range("B68").Value2 = 6361635
$ws.Range("F68").Value2 = "Aparecidense"
$ws.Range("G68").Value2 = "Sao Jose PA"
$ws.Range("H68").Value2 = 2
$ws.Range("I68").Value2 = 1
$ws.Range("J68").Value2 = "H"
$ws.Range("K68").Value2 = 2.15
$ws.Range("L68").Value2 = 3
$ws.Range("M68").Value2 = 3.2
$ws.Range("N68").Value2 = 2.3
$ws.Range("O68").Value2 = 3.1
$ws.Range("P68").Value2 = 3
$ws.Range("Q68").Value2 = -0.25
$ws.Range("R68").Value2 = 2
$ws.Range("S68").Value2 = 1.8
$ws.Range("T68").Value2 = 2
$ws.Range("U68").Value2 = 2.075
$ws.Range("V68").Value2 = 1.725
$ws.Range("W68").Value2 = 1.3
$ws.Range("X68").Value2 = -1
$ws.Range("Y68").Value2 = -1
$ws.Range("Z68").Value2 = 1
$ws.Range("AA68").Value2 = -1
$ws.Range("AB68").Value2 = 1.075
$ws.Range("AC68").Value2 = -1

# Row 69
$ws.Range("B69").Value2 = 6361463
$ws.Range("F69").Value2 = "AD Confianca"
$ws.Range("G69").Value2 = "Volta Redonda"
$ws.Range("H69").Value2 = 0
$ws.Range("I69").Value2 = 1
$ws.Range("J69").Value2 = "A"
$ws.Range("K69").Value2 = 2.3
$ws.Range("L69").Value2 = 2.875
$ws.Range("M69").Value2 = 3
$ws.Range("N69").Value2 = 2.3
$ws.Range("O69").Value2 = 3
$ws.Range("P69").Value2 = 3.1
$ws.Range("Q69").Value2 = -0.25
$ws.Range("R69").Value2 = 2.025
$ws.Range("S69").Value2 = 1.775
$ws.Range("T69").Value2 = 2
$ws.Range("U69").Value2 = 1.95
$ws.Range("V69").Value2 = 1.85
$ws.Range("W69").Value2 = -1
$ws.Range("X69").Value2 = -1
$ws.Range("Y69").Value2 = 2.1
$ws.Range("Z69").Value2 = -1
$ws.Range("AA69").Value2 = 0.7749999999999999
$ws.Range("AB69").Value2 = -1
$ws.Range("AC69").Value2 = 0.8500000000000001

# Row 70
$ws.Range("B70").Value2 = 6361634
$ws.Range("F70").Value2 = "Ypiranga"
$ws.Range("G70").Value2 = "Nautico"
$ws.Range("H70").Value2 = 3
$ws.Range("I70").Value2 = 3
$ws.Range("J70").Value2 = "D"
$ws.Range("K70").Value2 = 1.85
$ws.Range("L70").Value2 = 3.25
$ws.Range("M70").Value2 = 3.8
$ws.Range("N70").Value2 = 1.909
$ws.Range("O70").Value2 = 3.25
$ws.Range("P70").Value2 = 3.8
$ws.Range("Q70").Value2 = -0.5
$ws.Range("R70").Value2 = 1.975
$ws.Range("S70").Value2 = 1.825
$ws.Range("T70").Value2 = 2
$ws.Range("U70").Value2 = 1.975
$ws.Range("V70").Value2 = 1.825
$ws.Range("W70").Value2 = -1
$ws.Range("X70").Value2 = 2.25
$ws.Range("Y70").Value2 = -1
$ws.Range("Z70").Value2 = -1
$ws.Range("AA70").Value2 = 0.825
$ws.Range("AB70").Value2 = 0.9750000000000001
$ws.Range("AC70").Value2 = -1

# Row 82
$ws.Range("B82").Value2 = 6362496
$ws.Range("F82").Value2 = "Pouso Alegre"
$ws.Range("G82").Value2 = "Aparecidense"
$ws.Range("H82").Value2 = 0
$ws.Range("I82").Value2 = 1
$ws.Range("J82").Value2 = "A"
$ws.Range("K82").Value2 = 2.2
$ws.Range("L82").Value2 = 3
$ws.Range("M82").Value2 = 3.3
$ws.Range("N82").Value2 = 2.1
$ws.Range("O82").Value2 = 2.9
$ws.Range("P82").Value2 = 3.6
$ws.Range("Q82").Value2 = -0.25
$ws.Range("R82").Value2 = 1.875
$ws.Range("S82").Value2 = 1.925
$ws.Range("T82").Value2 = 1.75
$ws.Range("U82").Value2 = 1.775
$ws.Range("V82").Value2 = 2.025
$ws.Range("W82").Value2 = -1
$ws.Range("X82").Value2 = -1
$ws.Range("Y82").Value2 = 2.6
$ws.Range("Z82").Value2 = -1
$ws.Range("AA82").Value2 = 0.925
$ws.Range("AB82").Value2 = -1
$ws.Range("AC82").Value2 = 1.025

# Row 83
$ws.Range("B83").Value2 = 6361639
$ws.Range("F83").Value2 = "Amazonas FC"
$ws.Range("G83").Value2 = "Ypiranga"
$ws.Range("H83").Value2 = 0
$ws.Range("I83").Value2 = 0
$ws.Range("J83").Value2 = "D"
$ws.Range("K83").Value2 = 1.909
$ws.Range("L83").Value2 = 3.3
$ws.Range("M83").Value2 = 3.8
$ws.Range("N83").Value2 = 1.909
$ws.Range("O83").Value2 = 3.25
$ws.Range("P83").Value2 = 3.8
$ws.Range("Q83").Value2 = -0.5
$ws.Range("R83").Value2 = 1.95
$ws.Range("S83").Value2 = 1.85
$ws.Range("T83").Value2 = 2.25
$ws.Range("U83").Value2 = 2.025
$ws.Range("V83").Value2 = 1.775
$ws.Range("W83").Value2 = -1
$ws.Range("X83").Value2 = 2.25
$ws.Range("Y83").Value2 = -1
$ws.Range("Z83").Value2 = -1
$ws.Range("AA83").Value2 = 0.8500000000000001
$ws.Range("AB83").Value2 = -1
$ws.Range("AC83").Value2 = 0.7749999999999999

# Row 84
$ws.Range("B84").Value2 = 6361465
$ws.Range("F84").Value2 = "AD Confianca"
$ws.Range("G84").Value2 = "Operario PR"
$ws.Range("H84").Value2 = 1
$ws.Range("I84").Value2 = 0
$ws.Range("J84").Value2 = "H"
$ws.Range("K84").Value2 = 2.45
$ws.Range("L84").Value2 = 3
$ws.Range("M84").Value2 = 2.8
$ws.Range("N84").Value2 = 2.75
$ws.Range("O84").Value2 = 3
$ws.Range("P84").Value2 = 2.5
$ws.Range("Q84").Value2 = 0
$ws.Range("R84").Value2 = 2.025
$ws.Range("S84").Value2 = 1.775
$ws.Range("T84").Value2 = 2
$ws.Range("U84").Value2 = 1.95
$ws.Range("V84").Value2 = 1.85
$ws.Range("W84").Value2 = 1.75
$ws.Range("X84").Value2 = -1
$ws.Range("Y84").Value2 = -1
$ws.Range("Z84").Value2 = 1.025
$ws.Range("AA84").Value2 = -1
$ws.Range("AB84").Value2 = -1
$ws.Range("AC84").Value2 = 0.8500000000000001

# Row 90
$ws.Range("B90").Value2 = 6361948
$ws.Range("F90").Value2 = "Floresta EC"
$ws.Range("G90").Value2 = "Figueirense"
$ws.Range("H90").Value2 = 1
$ws.Range("I90").Value2 = 1
$ws.Range("J90").Value2 = "D"
$ws.Range("K90").Value2 = 2.55
$ws.Range("L90").Value2 = 3.1
$ws.Range("M90").Value2 = 2.625
$ws.Range("N90").Value2 = 2.1
$ws.Range("O90").Value2 = 3.1
$ws.Range("P90").Value2 = 3.4
$ws.Range("Q90").Value2 = -0.25
$ws.Range("R90").Value2 = 1.8
$ws.Range("S90").Value2 = 2
$ws.Range("T90").Value2 = 1.75
$ws.Range("U90").Value2 = 1.8
$ws.Range("V90").Value2 = 2
$ws.Range("W90").Value2 = -1
$ws.Range("X90").Value2 = 2.1
$ws.Range("Y90").Value2 = -1
$ws.Range("Z90").Value2 = -0.5
$ws.Range("AA90").Value2 = 0.5
$ws.Range("AB90").Value2 = 0.4
$ws.Range("AC90").Value2 = -0.5

# Row 91
$ws.Range("B91").Value2 = 6361638
$ws.Range("F91").Value2 = "Sao Jose PA"
$ws.Range("G91").Value2 = "Volta Redonda"
$ws.Range("H91").Value2 = 2
$ws.Range("I91").Value2 = 2
$ws.Range("J91").Value2 = "D"
$ws.Range("K91").Value2 = 2.2
$ws.Range("L91").Value2 = 3.2
$ws.Range("M91").Value2 = 3.1
$ws.Range("N91").Value2 = 2.1
$ws.Range("O91").Value2 = 3.2
$ws.Range("P91").Value2 = 3.3
$ws.Range("Q91").Value2 = -0.25
$ws.Range("R91").Value2 = 1.9
$ws.Range("S91").Value2 = 1.9
$ws.Range("T91").Value2 = 2
$ws.Range("U91").Value2 = 1.75
$ws.Range("V91").Value2 = 2.05
$ws.Range("W91").Value2 = -1
$ws.Range("X91").Value2 = 2.2
$ws.Range("Y91").Value2 = -1
$ws.Range("Z91").Value2 = -0.5
$ws.Range("AA91").Value2 = 0.45
$ws.Range("AB91").Value2 = 0.75
$ws.Range("AC91").Value2 = -1

# Row 104
$ws.Range("B104").Value2 = 6361641
$ws.Range("F104").Value2 = "Ypiranga"
$ws.Range("G104").Value2 = "Brusque"
$ws.Range("H104").Value2 = 2
$ws.Range("I104").Value2 = 0
$ws.Range("J104").Value2 = "H"
$ws.Range("K104").Value2 = 2.2
$ws.Range("L104").Value2 = 3.1
$ws.Range("M104").Value2 = 3
$ws.Range("N104").Value2 = 2.2
$ws.Range("O104").Value2 = 3
$ws.Range("P104").Value2 = 3
$ws.Range("Q104").Value2 = -0.25
$ws.Range("R104").Value2 = 2
$ws.Range("S104").Value2 = 1.8
$ws.Range("T104").Value2 = 2
$ws.Range("U104").Value2 = 1.9
$ws.Range("V104").Value2 = 1.9
$ws.Range("W104").Value2 = 1.2
$ws.Range("X104").Value2 = -1
$ws.Range("Y104").Value2 = -1
$ws.Range("Z104").Value2 = 1
$ws.Range("AA104").Value2 = -1
$ws.Range("AB104").Value2 = 0
$ws.Range("AC104").Value2 = -0

# Row 105
$ws.Range("B105").Value2 = 6361811
$ws.Range("F105").Value2 = "Clube Do Remo"
$ws.Range("G105").Value2 = "Volta Redonda"
$ws.Range("H105").Value2 = 2
$ws.Range("I105").Value2 = 1
$ws.Range("J105").Value2 = "H"
$ws.Range("K105").Value2 = 2.15
$ws.Range("L105").Value2 = 3.1
$ws.Range("M105").Value2 = 3.1
$ws.Range("N105").Value2 = 2.05
$ws.Range("O105").Value2 = 3
$ws.Range("P105").Value2 = 3.4
$ws.Range("Q105").Value2 = -0.25
$ws.Range("R105").Value2 = 1.8
$ws.Range("S105").Value2 = 2
$ws.Range("T105").Value2 = 2
$ws.Range("U105").Value2 = 1.825
$ws.Range("V105").Value2 = 1.975
$ws.Range("W105").Value2 = 1.05
$ws.Range("X105").Value2 = -1
$ws.Range("Y105").Value2 = -1
$ws.Range("Z105").Value2 = 0.8
$ws.Range("AA105").Value2 = -1
$ws.Range("AB105").Value2 = 0.825
$ws.Range("AC105").Value2 = -1

# Row 110
$ws.Range("B110").Value2 = 6361467
$ws.Range("F110").Value2 = "AD Confianca"
$ws.Range("G110").Value2 = "Floresta EC"
$ws.Range("H110").Value2 = 3
$ws.Range("I110").Value2 = 1
$ws.Range("J110").Value2 = "H"
$ws.Range("K110").Value2 = 1.95
$ws.Range("L110").Value2 = 3
$ws.Range("M110").Value2 = 3.75
$ws.Range("N110").Value2 = 1.95
$ws.Range("O110").Value2 = 3.1
$ws.Range("P110").Value2 = 3.6
$ws.Range("Q110").Value2 = -0.25
$ws.Range("R110").Value2 = 1.75
$ws.Range("S110").Value2 = 2.05
$ws.Range("T110").Value2 = 2
$ws.Range("U110").Value2 = 1.875
$ws.Range("V110").Value2 = 1.925
$ws.Range("W110").Value2 = 0.95
$ws.Range("X110").Value2 = -1
$ws.Range("Y110").Value2 = -1
$ws.Range("Z110").Value2 = 0.75
$ws.Range("AA110").Value2 = -1
$ws.Range("AB110").Value2 = 0.875
$ws.Range("AC110").Value2 = -1

# Row 111
$ws.Range("B111").Value2 = 6361642
$ws.Range("F111").Value2 = "Botafogo PB"
$ws.Range("G111").Value2 = "Sao Jose PA"
$ws.Range("H111").Value2 = 1
$ws.Range("I111").Value2 = 1
$ws.Range("J111").Value2 = "D"
$ws.Range("K111").Value2 = 1.833
$ws.Range("L111").Value2 = 3
$ws.Range("M111").Value2 = 4.333
$ws.Range("N111").Value2 = 1.85
$ws.Range("O111").Value2 = 3
$ws.Range("P111").Value2 = 4.333
$ws.Range("Q111").Value2 = -0.5
$ws.Range("R111").Value2 = 1.9
$ws.Range("S111").Value2 = 1.9
$ws.Range("T111").Value2 = 2
$ws.Range("U111").Value2 = 1.9
$ws.Range("V111").Value2 = 1.9
$ws.Range("W111").Value2 = -1
$ws.Range("X111").Value2 = 2
$ws.Range("Y111").Value2 = -1
$ws.Range("Z111").Value2 = -1
$ws.Range("AA111").Value2 = 0.8999999999999999
$ws.Range("AB111").Value2 = 0
$ws.Range("AC111").Value2 = -0

# Row 119
$ws.Range("B119").Value2 = 6361493
$ws.Range("F119").Value2 = "America RN"
$ws.Range("G119").Value2 = "Pouso Alegre"
$ws.Range("H119").Value2 = 2
$ws.Range("I119").Value2 = 0
$ws.Range("J119").Value2 = "H"
$ws.Range("K119").Value2 = 1.666
$ws.Range("L119").Value2 = 3.2
$ws.Range("M119").Value2 = 5
$ws.Range("N119").Value2 = 1.533
$ws.Range("O119").Value2 = 3.6
$ws.Range("P119").Value2 = 5.5
$ws.Range("Q119").Value2 = -1
$ws.Range("R119").Value2 = 2
$ws.Range("S119").Value2 = 1.8
$ws.Range("T119").Value2 = 2.25
$ws.Range("U119").Value2 = 1.95
$ws.Range("V119").Value2 = 1.85
$ws.Range("W119").Value2 = 0.5329999999999999
$ws.Range("X119").Value2 = -1
$ws.Range("Y119").Value2 = -1
$ws.Range("Z119").Value2 = 1
$ws.Range("AA119").Value2 = -1
$ws.Range("AB119").Value2 = -0.5
$ws.Range("AC119").Value2 = 0.425

# Row 120
$ws.Range("B120").Value2 = 6361643
$ws.Range("F120").Value2 = "Sao Jose PA"
$ws.Range("G120").Value2 = "Figueirense"
$ws.Range("H120").Value2 = 3
$ws.Range("I120").Value2 = 3
$ws.Range("J120").Value2 = "D"
$ws.Range("K120").Value2 = 2
$ws.Range("L120").Value2 = 3
$ws.Range("M120").Value2 = 3.5
$ws.Range("N120").Value2 = 1.85
$ws.Range("O120").Value2 = 3
$ws.Range("P120").Value2 = 4.2
$ws.Range("Q120").Value2 = -0.5
$ws.Range("R120").Value2 = 1.95
$ws.Range("S120").Value2 = 1.85
$ws.Range("T120").Value2 = 1.75
$ws.Range("U120").Value2 = 1.8
$ws.Range("V120").Value2 = 2
$ws.Range("W120").Value2 = -1
$ws.Range("X120").Value2 = 2
$ws.Range("Y120").Value2 = -1
$ws.Range("Z120").Value2 = -1
$ws.Range("AA120").Value2 = 0.8500000000000001
$ws.Range("AB120").Value2 = 0.8
$ws.Range("AC120").Value2 = -1

# Row 130
$ws.Range("B130").Value2 = 6361494
$ws.Range("F130").Value2 = "America RN"
$ws.Range("G130").Value2 = "Aparecidense"
$ws.Range("H130").Value2 = 0
$ws.Range("I130").Value2 = 0
$ws.Range("J130").Value2 = "D"
$ws.Range("K130").Value2 = 1.85
$ws.Range("L130").Value2 = 3.2
$ws.Range("M130").Value2 = 3.8
$ws.Range("N130").Value2 = 1.8
$ws.Range("O130").Value2 = 3
$ws.Range("P130").Value2 = 4.333
$ws.Range("Q130").Value2 = -0.5
$ws.Range("R130").Value2 = 1.875
$ws.Range("S130").Value2 = 1.925
$ws.Range("T130").Value2 = 1.75
$ws.Range("U130").Value2 = 1.75
$ws.Range("V130").Value2 = 2.05
$ws.Range("W130").Value2 = -1
$ws.Range("X130").Value2 = 2
$ws.Range("Y130").Value2 = -1
$ws.Range("Z130").Value2 = -1
$ws.Range("AA130").Value2 = 0.925
$ws.Range("AB130").Value2 = -1
$ws.Range("AC130").Value2 = 1.05

# Row 131
$ws.Range("B131").Value2 = 6361646
$ws.Range("F131").Value2 = "Sao Bernardo SP"
$ws.Range("G131").Value2 = "Ypiranga"
$ws.Range("H131").Value2 = 1
$ws.Range("I131").Value2 = 0
$ws.Range("J131").Value2 = "H"
$ws.Range("K131").Value2 = 1.909
$ws.Range("L131").Value2 = 3
$ws.Range("M131").Value2 = 3.9
$ws.Range("N131").Value2 = 1.7
$ws.Range("O131").Value2 = 3.1
$ws.Range("P131").Value2 = 5
$ws.Range("Q131").Value2 = -0.75
$ws.Range("R131").Value2 = 2
$ws.Range("S131").Value2 = 1.8
$ws.Range("T131").Value2 = 2
$ws.Range("U131").Value2 = 1.775
$ws.Range("V131").Value2 = 2.025
$ws.Range("W131").Value2 = 0.7
$ws.Range("X131").Value2 = -1
$ws.Range("Y131").Value2 = -1
$ws.Range("Z131").Value2 = 0.5
$ws.Range("AA131").Value2 = -0.5
$ws.Range("AB131").Value2 = -1
$ws.Range("AC131").Value2 = 1.025
